$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.388880320667681
$ws.Range("C2").Value = 0.26550508395664
$ws.Range("D2").Value = 0.02145684642999157
$ws.Range("F2").Value = 0.8335137090368221
$ws.Range("G2").Value = 0.6825438749194745
$ws.Range("H2").Value = 0.752941430975028
$ws.Range("I2").Value = 0.6861691456611112
$ws.Range("L2").Value = 0.3013015632718208
$ws.Range("B3").Value = 1.251468235563209
$ws.Range("C3").Value = 0.2427065503185588
$ws.Range("D3").Value = 0.02110144923565471
$ws.Range("F3").Value = 0.8252597966965141
$ws.Range("G3").Value = 0.6750462911114425
$ws.Range("H3").Value = 0.7556470654598684
$ws.Range("I3").Value = 0.6936146254535061
$ws.Range("L3").Value = 0.2899558975436491
$ws.Range("B4").Value = 1.167230185273581
$ws.Range("C4").Value = 0.2286129900187461
$ws.Range("D4").Value = 0.02088752658879756
$ws.Range("F4").Value = 0.8209591915412631
$ws.Range("G4").Value = 0.6711692591488685
$ws.Range("H4").Value = 0.757863823960605
$ws.Range("I4").Value = 0.6987925255703544
$ws.Range("L4").Value = 0.2831702641913552
$ws.Range("B5").Value = 1.13293725530508
$ws.Range("C5").Value = 0.2228459776758882
$ws.Range("D5").Value = 0.02080144232597547
$ws.Range("F5").Value = 0.8193987895404717
$ws.Range("G5").Value = 0.6697710005544195
$ws.Range("H5").Value = 0.7589065053231963
$ws.Range("I5").Value = 0.7010545096118932
$ws.Range("L5").Value = 0.2804504044245277
$ws.Range("B6").Value = 1.127245075305723
$ws.Range("C6").Value = 0.2218869386726396
$ws.Range("D6").Value = 0.02078721433820263
$ws.Range("F6").Value = 0.8191512660778599
$ws.Range("G6").Value = 0.6695497628045928
$ws.Range("H6").Value = 0.7590880470517192
$ws.Range("I6").Value = 0.7014392722923262
$ws.Range("L6").Value = 0.2800015098232791
$ws.Range("B7").Value = 1.16676755649604
$ws.Range("C7").Value = 0.2285353101087537
$ws.Range("D7").Value = 0.02088636119352927
$ws.Range("F7").Value = 0.8209373705826195
$ws.Range("G7").Value = 0.6711496675573443
$ws.Range("H7").Value = 0.75787732223273
$ws.Range("I7").Value = 0.6988224169458697
$ws.Range("L7").Value = 0.2831333997375367
$ws.Range("B8").Value = 1.341473578731325
$ws.Range("C8").Value = 0.2576639892168089
$ws.Range("D8").Value = 0.02133342116863091
$ws.Range("F8").Value = 0.8305080009408243
$ws.Range("G8").Value = 0.6798072990973765
$ws.Range("H8").Value = 0.7537588312873282
$ws.Range("I8").Value = 0.6886101854777849
$ws.Range("L8").Value = 0.2973520274271806
$ws.Range("B9").Value = 1.685097783977085
$ws.Range("C9").Value = 0.3140259539405008
$ws.Range("D9").Value = 0.02224373553696068
$ws.Range("F9").Value = 0.855403764358357
$ws.Range("G9").Value = 0.7025975905788471
$ws.Range("H9").Value = 0.7501065158696321
$ws.Range("I9").Value = 0.673418851984728
$ws.Range("L9").Value = 0.3266738463629366
$ws.Range("B10").Value = 1.938167138178528
$ws.Range("C10").Value = 0.3549707569273153
$ws.Range("D10").Value = 0.02293258260597852
$ws.Range("F10").Value = 0.8774881398252603
$ws.Range("G10").Value = 0.7229545887272479
$ws.Range("H10").Value = 0.7501441516080121
$ws.Range("I10").Value = 0.6652383202846366
$ws.Range("L10").Value = 0.3491046155645847
$ws.Range("B11").Value = 2.053425708417819
$ws.Range("C11").Value = 0.3734968017460574
$ws.Range("D11").Value = 0.02325022523000797
$ws.Range("F11").Value = 0.8883712524891649
$ws.Range("G11").Value = 0.7330151340279372
$ws.Range("H11").Value = 0.7507573682163979
$ws.Range("I11").Value = 0.6621711343178021
$ws.Range("L11").Value = 0.3595042783678792
$ws.Range("B12").Value = 2.097090023674866
$ws.Range("C12").Value = 0.3804976766184041
$ws.Range("D12").Value = 0.02337111610850684
$ws.Range("F12").Value = 0.8926136946955694
$ws.Range("G12").Value = 0.7369410111286356
$ws.Range("H12").Value = 0.7510757047997458
$ws.Range("I12").Value = 0.6611043229506421
$ws.Range("L12").Value = 0.3634706628824205
$ws.Range("B13").Value = 2.087685332222009
$ws.Range("C13").Value = 0.3789905620072886
$ws.Range("D13").Value = 0.02334505326293623
$ws.Range("F13").Value = 0.8916946020117251
$ws.Range("G13").Value = 0.7360903188263705
$ws.Range("H13").Value = 0.7510033084035115
$ws.Range("I13").Value = 0.6613298601087436
$ws.Range("L13").Value = 0.3626151725686952
$ws.Range("B14").Value = 2.057017635806972
$ws.Range("C14").Value = 0.3740730601821269
$ws.Range("D14").Value = 0.02326015889158839
$ws.Range("F14").Value = 0.8887178456432565
$ws.Range("G14").Value = 0.7333357837137839
$ws.Range("H14").Value = 0.7507818295547821
$ws.Range("I14").Value = 0.6620814664567902
$ws.Range("L14").Value = 0.3598300282805837
$ws.Range("B15").Value = 2.038235147828345
$ws.Range("C15").Value = 0.371059050508876
$ws.Range("D15").Value = 0.02320823735490052
$ws.Range("F15").Value = 0.8869103129431579
$ws.Range("G15").Value = 0.731663713087741
$ws.Range("H15").Value = 0.7506573954502613
$ws.Range("I15").Value = 0.6625541929556178
$ws.Range("L15").Value = 0.35812773009377
$ws.Range("B16").Value = 1.930637356713873
$ws.Range("C16").Value = 0.3537580109816645
$ws.Range("D16").Value = 0.02291190928185571
$ws.Range("F16").Value = 0.8767938195302634
$ws.Range("G16").Value = 0.722313303222009
$ws.Range("H16").Value = 0.7501161061864536
$ws.Range("I16").Value = 0.6654519827193184
$ws.Range("L16").Value = 0.3484289217697807
$ws.Range("B17").Value = 1.864663602253529
$ws.Range("C17").Value = 0.3431186659603895
$ws.Range("D17").Value = 0.02273121169369219
$ws.Range("F17").Value = 0.8708026551761634
$ws.Range("G17").Value = 0.7167828379252086
$ws.Range("H17").Value = 0.7499370063046342
$ws.Range("I17").Value = 0.6673976570735789
$ws.Range("L17").Value = 0.3425292370761497
$ws.Range("B18").Value = 1.826730117300087
$ws.Range("C18").Value = 0.3369897976633638
$ws.Range("D18").Value = 0.02262768284839822
$ws.Range("F18").Value = 0.8674354068715786
$ws.Range("G18").Value = 0.7136770974865527
$ws.Range("H18").Value = 0.7498900757727824
$ws.Range("I18").Value = 0.6685782963074232
$ws.Range("L18").Value = 0.3391543171963178
$ws.Range("B19").Value = 1.813888738247442
$ws.Range("C19").Value = 0.3349130590328571
$ws.Range("D19").Value = 0.02259269937539798
$ws.Range("F19").Value = 0.8663088054202177
$ws.Range("G19").Value = 0.7126384349895147
$ws.Range("H19").Value = 0.7498838058235862
$ws.Range("I19").Value = 0.6689885916831031
$ws.Range("L19").Value = 0.3380147883258928
$ws.Range("B20").Value = 1.871685297213503
$ws.Range("C20").Value = 0.3442522162933699
$ws.Range("D20").Value = 0.02275040555654684
$ws.Range("F20").Value = 0.8714322719904857
$ws.Range("G20").Value = 0.717363770554428
$ws.Range("H20").Value = 0.7499502644904652
$ws.Range("I20").Value = 0.667184163391461
$ws.Range("L20").Value = 0.3431553611601288
$ws.Range("B21").Value = 2.066024991527968
$ws.Range("C21").Value = 0.3755178454567272
$ws.Range("D21").Value = 0.02328507803797919
$ws.Range("F21").Value = 0.8895888933704299
$ws.Range("G21").Value = 0.7341416957670646
$ws.Range("H21").Value = 0.7508445426103663
$ws.Range("I21").Value = 0.6618581275132556
$ws.Range("L21").Value = 0.3606473249276547
$ws.Range("B22").Value = 2.193143922909371
$ws.Range("C22").Value = 0.3958669407365392
$ws.Range("D22").Value = 0.02363804865345998
$ws.Range("F22").Value = 0.9021624189332726
$ws.Range("G22").Value = 0.7457846708521458
$ws.Range("H22").Value = 0.7519311959102595
$ws.Range("I22").Value = 0.6589293143982289
$ws.Range("L22").Value = 0.3722440924697707
$ws.Range("B23").Value = 2.125288771044154
$ws.Range("C23").Value = 0.3850140576483057
$ws.Range("D23").Value = 0.02344934143415855
$ws.Range("F23").Value = 0.8953866888350035
$ws.Range("G23").Value = 0.739508221813594
$ws.Range("H23").Value = 0.7513051433037106
$ws.Range("I23").Value = 0.6604417667273808
$ws.Range("L23").Value = 0.3660395694150509
$ws.Range("B24").Value = 1.868510802213621
$ws.Range("C24").Value = 0.3437397760767453
$ws.Range("D24").Value = 0.02274172689926601
$ws.Range("F24").Value = 0.8711473819789717
$ws.Range("G24").Value = 0.7171009011277789
$ws.Range("H24").Value = 0.7499440959463755
$ws.Range("I24").Value = 0.6672804906655045
$ws.Range("L24").Value = 0.3428722378701394
$ws.Range("B25").Value = 1.592030339315158
$ws.Range("C25").Value = 0.2988600037775768
$ws.Range("D25").Value = 0.02199392076740736
$ws.Range("F25").Value = 0.8480068078429781
$ws.Range("G25").Value = 0.6958027524671877
$ws.Range("H25").Value = 0.7506184660874453
$ws.Range("I25").Value = 0.6770074897044509
$ws.Range("L25").Value = 0.3185863503169202
